$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 07.01.2022 15:18"

# D2: change from text "+0.6" to numeric value 0.6
$ws.Range("D2").Value = 0.6

# E2: change from text date to numeric Excel date serial value,
# using the same date/time number format style as the other rows (E3 etc.)
$ws.Range("E2").Value = 44568.636875
$ws.Range("E2").NumberFormat = $ws.Range("E3").NumberFormat
